# Generate Report for Handoff
# Mark the f8a1579f-8434-4b75-924c-95fa65df49b9 file as "Ready for handoff"
# across the Overview, zh-cn and de-de sheets, and stamp the new handoff
# datetimes.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is f8a1579f-8434-4b75-924c-95fa65df49b9.md ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"
$overview.Range("D3").Value = "2016-17-21 00:17:27"

# --- zh-cn sheet: row 3 is f8a1579f-8434-4b75-924c-95fa65df49b9 ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "2016-03-21 00:17:24"

# --- de-de sheet: row 3 is f8a1579f-8434-4b75-924c-95fa65df49b9 ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "2016-03-21 00:17:27"
